# Update sheet for 2023-11 (row 43): fill in the monthly figures for
# November 2023 (fuel, toll, parking, violation, maintenance, rideshare
# income, fuel refund) and extend the "difference" shared formula down
# into the newly filled row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data for row 43 (2023 / month 11)
$ws.Range("C43").Value = 1027.29
$ws.Range("D43").Value = 593
$ws.Range("E43").Value = 55
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0

# Extend the J column "difference" formula from the last filled row (J34:J42)
# down through the new row (J43), matching J34's pattern: (H+I)-(C+D+E+F+G)
$ws.Range("J34:J43").Formula = "=(H34+I34)-(C34+D34+E34+F34+G34)"

# Reflect the active cell / selection on the newly completed row
$ws.Range("J43").Select()
